$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (e.g. "537.18") need to be
# forced to Text first, otherwise Excel auto-converts them to real numbers (losing
# exact formatting like trailing zeros) instead of keeping the original inline-string type.
# NumberFormat is reset via ClearFormats() right after the write so no residual style
# ends up attached to the cell (matches the source file, which carries no style on these cells).
$textForceCells = @("D5", "D19", "D50", "D27", "D37", "D31", "D33", "D51", "D18", "D15", "D21", "D42", "D6", "D17", "D44", "D47", "D40", "D20", "D36", "D9", "D26", "D23", "D32")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '59.540.52'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.645.49'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '537.18'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '145.68'
$ws.Range("E6").Value = '  +3.51%  '
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("D9").Value = '6.69'
$ws.Range("E9").Value = '  +3.68%  '
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("D13").Value = '3.109.27'
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = '59.459.25'
$ws.Range("E14").Value = '  +0.32%  '
$ws.Range("D15").Value = '21.24'
$ws.Range("E15").Value = '  +3.51%  '
$ws.Range("D16").Value = '2.646.14'
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").Value = '339.28'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").Value = '4.40'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = '10.31'
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").Value = '6.27'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '66.76'
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '7.27'
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").Value = '0.0₃0744'
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = '5.85'
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").Value = '18.91'
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = '151.24'
$ws.Range("E33").Value = '  +1.33%  '
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").Value = '0.846'
$ws.Range("E36").Value = '  +2.41%  '
$ws.Range("D37").Value = '0.834'
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("E39").Value = '  +1.54%  '
$ws.Range("D40").Value = '286.75'
$ws.Range("E40").Value = '  +4.49%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '0.606'
$ws.Range("E42").Value = '  +1.79%  '
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").Value = '0.0540'
$ws.Range("E44").Value = '  +2.99%  '
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0226'
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.967.55'
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("D50").Value = '18.27'
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").Value = '110.90'
$ws.Range("E51").Value = '  -0.04%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
